$d = $word.ActiveDocument

# Find the "Marketing Strategy and Data-Driven Insights" paragraph under the
# Siege Analytics PARTNER role, then insert three new bullet paragraphs right
# after it (before the existing "Conducted comprehensive..." bullet).
$anchor = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    $t = $t -replace "[\r\a]+$", ""
    if ($t -eq "Marketing Strategy and Data-Driven Insights") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Could not find 'Marketing Strategy and Data-Driven Insights' paragraph"
}

$r = $anchor.Range
$r.Collapse(0)  # wdCollapseEnd - collapse to the end of the paragraph (after its mark)

$newLines = @(
    "• Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters",
    "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
    "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
)

for ($i = $newLines.Length - 1; $i -ge 0; $i--) {
    $r.InsertBefore($newLines[$i] + "`r")
}
